# Insert a new "Domingo" day column into the weekly schedule sheet.
#
# Net effect (per the target OOXML):
#   - A brand-new column H is appended (it inherits column G's row-1
#     header style), extending the sheet's used range by one column.
#   - Row 1's headers are rewritten in place so "Domingo" becomes the new
#     B1, and the remaining weekday names shift one column to the right
#     (Jueves->C1, Lunes->D1, Martes->E1, Miercoles->F1, Sabado->G1,
#     Viernes->H1).
#   - Every data row (2-55) gets an explicit (empty) H cell, matching the
#     already-explicit empty B:G cells those rows already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 55

# 1) Make room for the new trailing column (H) by inserting an entire
#    column at H; since H is one past the previous last column (G) this
#    just materialises the new column rather than displacing data.
$ws.Columns.Item(8).Insert()

# 2) Rewrite the header row with "Domingo" inserted after "Salon" and the
#    rest of the week shifted right by one column.
$headers = @("Salon", "Domingo", "Jueves", "Lunes", "Martes", "Miercoles", "Sabado", "Viernes")
for ($col = 1; $col -le 8; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# 3) Materialise an empty (but present) text cell in column H on every
#    data row, matching the already-present empty B:G cells on those
#    rows. A lone leading apostrophe is Excel's classic "force text,
#    empty content" entry, then formatting is reset so no stray style is
#    left applied to the cell.
for ($row = 2; $row -le $lastRow; $row++) {
    $c = $ws.Cells.Item($row, 8)
    $c.Formula = "'"
    $c.ClearFormats()
}
